$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before N (shifts old N/O/P -> O/P/Q) and give it the
# same visual width as its left neighbour (column M), matching Excel's
# "insert column" behaviour.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 10.2

# Make "Repayment schedule" the active sheet/tab with the new selection.
$ws.Range("L12").Select()
